$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Remove all existing hyperlinks first (clean slate) so that re-adding
#    them below produces a deterministic, correctly-ordered relationship list.
# ---------------------------------------------------------------------------
$ws.Cells.Item(3, 2).Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 2) Update/insert the endpoint table rows (A=METODO, B=ENDPOINT, C=DESCRIPCION)
#    Base URL moved from the old public IP to the local dev server.
# ---------------------------------------------------------------------------

# Row 3 - GET swagger-ui (existing row, URL changed)
$ws.Range("A3").Value2 = "GET"
$ws.Range("B3").Value2 = "http://localhost:9090/swagger-ui.html"
$ws.Range("C3").Value2 = "Documentación de la API REST"

# Row 4 - POST usuarios (existing row, URL changed)
$ws.Range("A4").Value2 = "POST"
$ws.Range("B4").Value2 = "http://localhost:9090/usuarios"
$ws.Range("C4").Value2 = "Servicio para registrar usuarios"

# Row 5 - POST oauth/token (existing row, URL changed)
$ws.Range("A5").Value2 = "POST"
$ws.Range("B5").Value2 = "http://localhost:9090/oauth/token"
$ws.Range("C5").Value2 = "Servicio para inicio de sesion (Oauth JWT)"

# Row 6 - POST vehiculos (new service)
$ws.Range("A6").Value2 = "POST"
$ws.Range("B6").Value2 = "http://localhost:9090/vehiculos"
$ws.Range("C6").Value2 = "Servicio para registro de vehiculos"

# Row 7 - POST departamentos (new service)
$ws.Range("A7").Value2 = "POST"
$ws.Range("B7").Value2 = "http://localhost:9090/departamentos"
$ws.Range("C7").Value2 = "Servicio para registro de departamentos"

# Row 8 - GET departamentos (new service)
$ws.Range("A8").Value2 = "GET"
$ws.Range("B8").Value2 = "http://localhost:9090/departamentos"
$ws.Range("C8").Value2 = "Servicio para consulta de departamentos"

# Row 9 - POST conductores (new service, highlighted)
$ws.Range("A9").Value2 = "POST"
$ws.Range("B9").Value2 = "http://localhost:9090/conductores"
$ws.Range("C9").Value2 = "Servicio para registro de conductores"

# Row 10 - GET conductores (new service, highlighted)
$ws.Range("A10").Value2 = "GET"
$ws.Range("B10").Value2 = "http://localhost:9090/conductores"
$ws.Range("C10").Value2 = "Servicio para consulta de conductores"

# Row 11 - POST encomiendas (new service, highlighted)
$ws.Range("A11").Value2 = "POST"
$ws.Range("B11").Value2 = "http://localhost:9090/encomiendas"
$ws.Range("C11").Value2 = "Servicio para registro de encomiendas"

# Row 12 - POST rutas (new service, highlighted)
$ws.Range("A12").Value2 = "POST"
$ws.Range("B12").Value2 = "http://localhost:9090/rutas"
$ws.Range("C12").Value2 = "Servicio para registro de programacion de rutas"

# ---------------------------------------------------------------------------
# 3) Re-create the hyperlinks in the same order as the original authoring so
#    relationship ids line up: B4, B5, B3 (existing three), then the seven
#    newly-added ones in insertion order.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("B4"), "http://localhost:9090/usuarios") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B5"), "http://localhost:9090/oauth/token") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "http://localhost:9090/swagger-ui.html") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B6"), "http://localhost:9090/vehiculos") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B7"), "http://localhost:9090/departamentos") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B9"), "http://localhost:9090/conductores") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B8"), "http://localhost:9090/departamentos") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B10"), "http://localhost:9090/conductores") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B11"), "http://localhost:9090/encomiendas") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B12"), "http://localhost:9090/rutas") | Out-Null

# ---------------------------------------------------------------------------
# 4) Highlight the newest block of services (conductores / encomiendas /
#    rutas - rows 9-12) with a yellow fill, matching the existing bordered
#    look of the rest of the table.
# ---------------------------------------------------------------------------
$ws.Range("A9:C12").Interior.Color = 65535

# ---------------------------------------------------------------------------
# 5) Selection, matching the author's last recorded cursor position.
# ---------------------------------------------------------------------------
$ws.Range("B15").Select() | Out-Null
